# Daily attendance processing - reverse order of names in the
# "Recorded By" column (G) for every data row that lists multiple
# recorders, e.g. "a, b, c" becomes "c, b, a".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -is [string] -and $value.Contains(",")) {
        $parts = $value -split ",\s*"
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }
        $reversed = $trimmed[($trimmed.Length - 1)..0]
        $newValue = [string]::Join(", ", $reversed)
        $cell.Value = $newValue
    }
}
